# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
# Swap the data (columns B..AC) between row 11 and row 12, and between row 83
# and row 84, leaving column A (the row's serial id) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($Sheet, $RowA, $RowB)

    # Columns B (2) through AC (29) hold the data that must be swapped.
    $firstCol = 2
    $lastCol = 29

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $Sheet.Cells.Item($RowA, $col)
        $cellB = $Sheet.Cells.Item($RowB, $col)

        $valueA = $cellA.Value2
        $valueB = $cellB.Value2

        $cellA.Value = $valueB
        $cellB.Value = $valueA
    }
}

Swap-RowData $ws 11 12
Swap-RowData $ws 83 84
